# Insert a new data row above the current row 132 ("Fruta, Macroferia
# Regional de Talca - Arándano (blue)" weekly price update). All existing
# rows from 132 downward shift down by one (132->133, ..., 159->160).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 132 (and everything below it) down by one row, carrying
# formatting along (matches Excel's native "insert row" behaviour, which is
# what happened in the source edit: the table grew from 159 to 160 rows).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new record.
$ws.Cells.Item(132, 1).Value  = 5
$ws.Cells.Item(132, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(132, 3).Value  = "Maule"
$ws.Cells.Item(132, 4).Value  = 45275
$ws.Cells.Item(132, 5).Value  = 7
$ws.Cells.Item(132, 6).Value  = "Fruta"
$ws.Cells.Item(132, 7).Value  = 100101
$ws.Cells.Item(132, 8).Value  = "Berries"
$ws.Cells.Item(132, 9).Value  = 100101001
$ws.Cells.Item(132, 10).Value = "Arándano (blue)"
$ws.Cells.Item(132, 11).Value = "Sin especificar"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 120
$ws.Cells.Item(132, 14).Value = 4000
$ws.Cells.Item(132, 15).Value = 4000
$ws.Cells.Item(132, 16).Value = 4000
$ws.Cells.Item(132, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(132, 18).Value = "Provincia de Linares"
$ws.Cells.Item(132, 19).Value = 2000
$ws.Cells.Item(132, 20).Value = 2
